$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bets")

# M2 used to show the raw win/loss % change (F2/D2-1) formatted as 0.000%.
# Re-express the same ratio as a plain number scaled to "percentage points"
# (multiply by 100, rounded to 3 decimals) and drop the custom percent
# number format in favour of the default/general display.
$ws.Range("M2").NumberFormat = "General"
$ws.Range("M2").Formula = "=ROUND((F2/D2-1)*100, 3)"

# Cursor ended up on J8 after the edit.
$ws.Range("J8").Select() | Out-Null
